$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.426.60'
$ws.Range("E2").Value = '  +0.61%  '

$ws.Range("D3").Value = '2.648.07'
$ws.Range("E3").Value = '  +0.41%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '598.24'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.03%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.88'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.56%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("E8").Value = '  -0.83%  '

$ws.Range("D9").Value = '2.646.54'
$ws.Range("E9").Value = '  +0.43%  '

$ws.Range("E10").Value = '  +7.73%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.28'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.02%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.356'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.18%  '

$ws.Range("E14").Value = '  +1.62%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000193'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.72%  '

$ws.Range("D16").Value = '3.129.45'
$ws.Range("E16").Value = '  +0.40%  '

$ws.Range("D17").Value = '68.270.47'
$ws.Range("E17").Value = '  +0.60%  '

$ws.Range("D18").Value = '2.636.52'
$ws.Range("E18").Value = '  -0.35%  '

$ws.Range("E19").Value = '  -0.56%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '364.22'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.18%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.50'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.09%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.39'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.31%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.91'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.99%  '

$ws.Range("E24").Value = '  +0.36%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '74.74'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.49%  '

$ws.Range("E26").Value = '  +0.00%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.85'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.07%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0000107'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.31%  '

$ws.Range("D29").Value = '2.776.39'
$ws.Range("E29").Value = '  +0.36%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.23%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '574.57'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.02%  '

$ws.Range("E32").Value = '  +2.85%  '

$ws.Range("E33").Value = '  +1.72%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.89'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.26%  '

$ws.Range("E35").Value = '  +3.11%  '

$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.60'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.71%  '

$ws.Range("B37").Value = 'FirstDigitalUSD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.03%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '160.99'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.84%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.37'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.96%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.376'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.98%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.90'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.26%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.39'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.36%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.67'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.15%  '

$ws.Range("E44").Value = '  -0.69%  '

$ws.Range("E45").Value = '  +2.52%  '

$ws.Range("E46").Value = '  +0.00%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '40.61'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.85%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '157.01'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.31%  '

$ws.Range("E49").Value = '  +1.82%  '

$ws.Range("E50").Value = '  +0.89%  '

$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '21.92'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.34%  '
